$d = $word.ActiveDocument

$d.Content.Find.Execute("47×87=", $true, $false, $false, $false, $false, $true, 1, $false, "64×42=", 2) | Out-Null
$d.Content.Find.Execute("43×86=", $true, $false, $false, $false, $false, $true, 1, $false, "71×82=", 2) | Out-Null
$d.Content.Find.Execute("77×83=", $true, $false, $false, $false, $false, $true, 1, $false, "81×80=", 2) | Out-Null
$d.Content.Find.Execute("23×29=", $true, $false, $false, $false, $false, $true, 1, $false, "87×49=", 2) | Out-Null
$d.Content.Find.Execute("55×81=", $true, $false, $false, $false, $false, $true, 1, $false, "39×77=", 2) | Out-Null
$d.Content.Find.Execute("79×26=", $true, $false, $false, $false, $false, $true, 1, $false, "79×79=", 2) | Out-Null
$d.Content.Find.Execute("36×55=", $true, $false, $false, $false, $false, $true, 1, $false, "86×89=", 2) | Out-Null
$d.Content.Find.Execute("97×14=", $true, $false, $false, $false, $false, $true, 1, $false, "38×96=", 2) | Out-Null
$d.Content.Find.Execute("73×19=", $true, $false, $false, $false, $false, $true, 1, $false, "68×66=", 2) | Out-Null
$d.Content.Find.Execute("71×95=", $true, $false, $false, $false, $false, $true, 1, $false, "69×89=", 2) | Out-Null
$d.Content.Find.Execute("84×11=", $true, $false, $false, $false, $false, $true, 1, $false, "76×75=", 2) | Out-Null
$d.Content.Find.Execute("24×68=", $true, $false, $false, $false, $false, $true, 1, $false, "43×46=", 2) | Out-Null
$d.Content.Find.Execute("43×43=", $true, $false, $false, $false, $false, $true, 1, $false, "41×55=", 2) | Out-Null
$d.Content.Find.Execute("65×23=", $true, $false, $false, $false, $false, $true, 1, $false, "47×16=", 2) | Out-Null
$d.Content.Find.Execute("15×56=", $true, $false, $false, $false, $false, $true, 1, $false, "57×73=", 2) | Out-Null
$d.Content.Find.Execute("45×23=", $true, $false, $false, $false, $false, $true, 1, $false, "63×57=", 2) | Out-Null
$d.Content.Find.Execute("61×55=", $true, $false, $false, $false, $false, $true, 1, $false, "19×77=", 2) | Out-Null
$d.Content.Find.Execute("29×47=", $true, $false, $false, $false, $false, $true, 1, $false, "82×65=", 2) | Out-Null
$d.Content.Find.Execute("34×41=", $true, $false, $false, $false, $false, $true, 1, $false, "96×28=", 2) | Out-Null
$d.Content.Find.Execute("88×78=", $true, $false, $false, $false, $false, $true, 1, $false, "45×44=", 2) | Out-Null
$d.Content.Find.Execute("53×95=", $true, $false, $false, $false, $false, $true, 1, $false, "44×45=", 2) | Out-Null
$d.Content.Find.Execute("28×78=", $true, $false, $false, $false, $false, $true, 1, $false, "68×28=", 2) | Out-Null
$d.Content.Find.Execute("95×42=", $true, $false, $false, $false, $false, $true, 1, $false, "56×64=", 2) | Out-Null
$d.Content.Find.Execute("84×13=", $true, $false, $false, $false, $false, $true, 1, $false, "59×86=", 2) | Out-Null
$d.Content.Find.Execute("63×44=", $true, $false, $false, $false, $false, $true, 1, $false, "28×12=", 2) | Out-Null
